$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.617.05"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.167.82"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.73"
$ws.Range("D5").NumberFormat = $null
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.07"
$ws.Range("D7").NumberFormat = $null
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("D9").NumberFormat = $null
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0850"
$ws.Range("D10").NumberFormat = $null
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "2.490.16"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.72"
$ws.Range("D14").NumberFormat = $null
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.810"
$ws.Range("D15").NumberFormat = $null
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").NumberFormat = $null
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "2.168.93"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "39.565.65"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "0.0₃0911"
$ws.Range("E19").Value = "  +6.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.70"
$ws.Range("D20").NumberFormat = $null
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.72"
$ws.Range("D22").NumberFormat = $null
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("D24").NumberFormat = $null
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.81"
$ws.Range("D26").NumberFormat = $null
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("E27").Value = "  -3.25%  "
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.69"
$ws.Range("D30").NumberFormat = $null
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("D31").NumberFormat = $null
$ws.Range("E31").Value = "  +3.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("D33").NumberFormat = $null
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.96"
$ws.Range("D35").NumberFormat = $null
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.80"
$ws.Range("D37").NumberFormat = $null
$ws.Range("E37").Value = "  +6.16%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").NumberFormat = $null
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.95"
$ws.Range("D40").NumberFormat = $null
$ws.Range("E40").Value = "  +18.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.50"
$ws.Range("D41").NumberFormat = $null
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.66"
$ws.Range("D43").NumberFormat = $null
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("D44").Value = "1.514.03"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.89"
$ws.Range("D46").NumberFormat = $null
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000195"
$ws.Range("D50").NumberFormat = $null
$ws.Range("E50").Value = "  +31.86%  "
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.99"
$ws.Range("D51").NumberFormat = $null
$ws.Range("E51").Value = "  +0.54%  "
